# "New query for Union" — add a new SQL-query row (row 11) to the
# "10 SQL Retrieval Queries" checklist sheet, tick its "Completed"
# checkbox, and move the sheet selection onto the new description cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Row 11: mark the checkbox cell TRUE (checkbox "Check Box 6" is
# form-linked to $A$11) and add the new query description in B11.
$ws.Range("A11").Value = $true
$ws.Range("B11").Value = "Heaviest and lightest Pokemon"

# Reflect the checked state on the checkbox shape itself too.
try {
    $ws.Shapes.Item("Check Box 6").ControlFormat.Value = 1
} catch {
}

# Move the active selection to the newly-filled cell, like the user
# would after typing the new query text.
$null = $ws.Range("B11").Select()
